# Update cryptos list values per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$sub3 = [char]0x2083

$ws.Range("D2").Value = "41.799.06"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "2.478.61"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.96"
$ws.Range("E5").Value = "  +1.53%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.48"
$ws.Range("E6").Value = "  +2.19%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.554"
$ws.Range("E7").Value = "  +0.74%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  +1.28%  "

$ws.Range("E10").Value = "  +11.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.32"
$ws.Range("E11").Value = "  +2.41%  "

$ws.Range("E12").Value = "  +0.62%  "

$ws.Range("D13").Value = "2.859.22"
$ws.Range("E13").Value = "  +0.53%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.94"
$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.73"
$ws.Range("E15").Value = "  -0.43%  "

$ws.Range("D16").Value = "2.467.01"
$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.805"
$ws.Range("E17").Value = "  +3.85%  "

$ws.Range("D18").Value = "41.754.68"

$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.50"
$ws.Range("E19").Value = "  +0.04%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0$($sub3)0955"
$ws.Range("E20").Value = "  +1.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.33"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.37"
$ws.Range("E22").Value = "  +2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.02"
$ws.Range("E23").Value = "  +1.65%  "

$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("E25").Value = "  +2.96%  "

$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.31"
$ws.Range("E27").Value = "  +3.23%  "

$ws.Range("E28").Value = "  +0.73%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.77"
$ws.Range("E29").Value = "  +1.10%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.13"
$ws.Range("E30").Value = "  +5.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.99"
$ws.Range("E31").Value = "  +1.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.54"
$ws.Range("E32").Value = "  +1.92%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0767"
$ws.Range("E34").Value = "  +1.21%  "

$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.47"
$ws.Range("E36").Value = "  +1.49%  "

$ws.Range("E37").Value = "  +5.36%  "

$ws.Range("E38").Value = "  +1.90%  "

$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("E40").Value = "  +1.63%  "

$ws.Range("E41").Value = "  +1.28%  "

$ws.Range("E42").Value = "  +8.05%  "

$ws.Range("D43").Value = "2.009.31"
$ws.Range("E43").Value = "  +3.44%  "

$ws.Range("E44").Value = "  +3.72%  "

$ws.Range("E45").Value = "  +1.02%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.99"
$ws.Range("E46").Value = "  +3.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.51"
$ws.Range("E47").Value = "  +5.54%  "

$ws.Range("D48").Value = "2.716.47"
$ws.Range("E48").Value = "  +0.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.98"
$ws.Range("E49").Value = "  +7.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.45"
$ws.Range("E50").Value = "  +1.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.49"
$ws.Range("E51").Value = "  +0.87%  "
